$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.415.04'
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').Value = '2.578.66'
$ws.Range('E3').Value = '  -2.65%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.26'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '150.53'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.25%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.587'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.72'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.56'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = '3.042.35'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '63.252.81'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000156'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +5.71%  '
$ws.Range('D17').Value = '2.580.81'
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '346.04'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.85'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.83%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.21'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.70'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.67'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.28%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.15'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '554.32'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.07'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.162'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.04'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').Value = '0.0₃0862'
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.24'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '166.77'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.54'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '165.67'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -1.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.98'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '22.95'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.91%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0586'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.94%  '
$ws.Range('E46').Value = '  +5.57%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.627'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0963'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('D51').Value = '0.0₆0234'
$ws.Range('E51').Value = '  +18.17%  '
